# Automatische test-sync: 2025-06-27 22:42:50
# Append the new "Wanneer zijn jullie open?" test-mail row (row 10) to the
# "Logs" sheet, extend the conditional-formatting ranges to cover it, and
# bump the "Openingstijden / Locatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 10

$ws.Cells.Item($row, 1).Value = "Wanneer zijn jullie open?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Cells.Item($row, 4).Value = "Openingstijden / Locatie"
$ws.Cells.Item($row, 5).Value = "Beste klant,`n`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`n`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Cells.Item($row, 6).Value = "2025-06-27 22:42:34"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"

# Extend the conditional-formatting ranges (columns D, G, H, I) from row 9
# down to the freshly added row 10, mirroring each rule's new sqref.
$dFcs = $ws.Range("D2:D9").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D10"))
}

$gFcs = $ws.Range("G2:G9").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G10"))
}

$hFcs = $ws.Range("H2:H9").FormatConditions
for ($i = 1; $i -le $hFcs.Count; $i++) {
    $hFcs.Item($i).ModifyAppliesToRange($ws.Range("H2:H10"))
}

$iFcs = $ws.Range("I2:I9").FormatConditions
for ($i = 1; $i -le $iFcs.Count; $i++) {
    $iFcs.Item($i).ModifyAppliesToRange($ws.Range("I2:I10"))
}

# Dashboard: "Openingstijden / Locatie" count goes from 6 to 7.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 7
